$wb = $excel.ActiveWorkbook

# Data for 2024-11-15: updated YTD crime counts across worksheets
$updates = @{
    'Citywide Totals' = @{ 'F2'=82; 'H2'=100; 'I2'=112; 'E3'=137; 'H3'=143; 'J3'=212; 'K3'=206; 'E4'=11; 'J4'=20; 'D6'=387; 'E6'=437; 'F6'=485; 'G6'=420; 'I6'=477; 'J6'=389; 'K6'=470; 'D7'=607; 'E7'=654; 'F7'=700; 'G7'=640; 'H7'=682; 'I7'=798; 'J7'=734; 'K7'=834 }
    'Garfield Park' = @{ 'K3'=16; 'E6'=48; 'E7'=61; 'K7'=44 }
    'Grand Crossing' = @{ 'E3'=8; 'K6'=43; 'E7'=34; 'K7'=64 }
    'Washington Park' = @{ 'E4'=4; 'E5'=8 }
    'By Neighborhood' = @{ 'E2'=4; 'E6'=2; 'E8'=47; 'F8'=43; 'G8'=32; 'J8'=43; 'F19'=22; 'G29'=10; 'K29'=19; 'E32'=61; 'K32'=44; 'E36'=34; 'K36'=64; 'E47'=16; 'E48'=6; 'D53'=68; 'E53'=80; 'H53'=92; 'I53'=122; 'J65'=11; 'J70'=12; 'J74'=21; 'D81'=3; 'E88'=8; 'D98'=607; 'E98'=654; 'F98'=700; 'G98'=640; 'H98'=682; 'I98'=798; 'J98'=734; 'K98'=834 }
    'Loop' = @{ 'H2'=12; 'I2'=12; 'H3'=18; 'D6'=40; 'E6'=62; 'I6'=78; 'D7'=68; 'E7'=80; 'H7'=92; 'I7'=122 }
    'South Deering' = @{ 'D5'=2; 'D6'=3 }
    'Rogers Park' = @{ 'J4'=1 }
    'North Lawndale' = @{ 'J6'=11 }
    'River North' = @{ 'J3'=5; 'J6'=21 }
    'Fuller Park' = @{ 'G5'=8; 'K5'=15; 'G6'=10; 'K6'=19 }
    'Chatham' = @{ 'F6'=15; 'F7'=22 }
    'Albany Park' = @{ 'E4'=2; 'E6'=4 }
    'Lincoln Park' = @{ 'E5'=5; 'E6'=6 }
    'Lake View' = @{ 'E5'=10; 'E6'=16 }
    'Old Town' = @{ 'J5'=5; 'J6'=12 }
    'Austin' = @{ 'F2'=8; 'E6'=37; 'G6'=23; 'J6'=23; 'E7'=47; 'F7'=43; 'G7'=32; 'J7'=43 }
    'Ashburn' = @{ 'E4'=2; 'E5'=2 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($addr in $cellMap.Keys) {
        $ws.Range($addr).Value = $cellMap[$addr]
    }
}
